# Update "想去人数" (column F) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 122
    4  = 1632
    5  = 626
    6  = 1100
    8  = 11566
    11 = 454
    12 = 370
    13 = 1091
    14 = 811
    15 = 12401
    16 = 13144
    21 = 238
    24 = 127
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
